# Insert a new pricing record as row 83, pushing the existing rows
# 83..181 down to 84..182 (weekly fruit/vegetable price update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 83.
$ws.Rows.Item(83).Insert()

# Populate the new row 83 with the new observation.
$ws.Range("A83").Value = 7
$ws.Range("B83").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C83").Value = "Ñuble"
$ws.Range("D83").Value = 45159
$ws.Range("E83").Value = 16
$ws.Range("F83").Value = "Fruta"
$ws.Range("G83").Value = 100108
$ws.Range("H83").Value = "Tropicales y subtropicales"
$ws.Range("I83").Value = 100108002
$ws.Range("J83").Value = "Mango"
$ws.Range("K83").Value = "Sin especificar"
$ws.Range("L83").Value = "Primera"
$ws.Range("M83").Value = 50
$ws.Range("N83").Value = 8000
$ws.Range("O83").Value = 8000
$ws.Range("P83").Value = 8000
$ws.Range("Q83").Value = "$/bandeja 4 kilos"
$ws.Range("R83").Value = "Brasil"
$ws.Range("S83").Value = 2000
$ws.Range("T83").Value = 4
